$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2024-07-07 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-08 Monday", 2) | Out-Null

# Update table cell contents (addressed by row/col to disambiguate duplicate values)
$t = $d.Tables(1)
$t.Cell(1, 1).Range.Text = "44÷2=22, 0"
$t.Cell(1, 2).Range.Text = "79÷2=39, 1"
$t.Cell(1, 3).Range.Text = "12÷6=2, 0"
$t.Cell(1, 4).Range.Text = "47÷2=23, 1"
$t.Cell(1, 5).Range.Text = "47÷6=7, 5"
$t.Cell(5, 1).Range.Text = "63÷2=31, 1"
$t.Cell(5, 2).Range.Text = "39÷5=7, 4"
$t.Cell(5, 3).Range.Text = "64÷7=9, 1"
$t.Cell(5, 4).Range.Text = "33÷2=16, 1"
$t.Cell(5, 5).Range.Text = "31÷6=5, 1"
$t.Cell(9, 1).Range.Text = "10÷6=1, 4"
$t.Cell(9, 2).Range.Text = "66÷2=33, 0"
$t.Cell(9, 3).Range.Text = "81÷7=11, 4"
$t.Cell(9, 4).Range.Text = "23÷7=3, 2"
$t.Cell(9, 5).Range.Text = "63÷2=31, 1"
$t.Cell(13, 1).Range.Text = "92÷7=13, 1"
$t.Cell(13, 2).Range.Text = "79÷9=8, 7"
$t.Cell(13, 3).Range.Text = "52÷2=26, 0"
$t.Cell(13, 4).Range.Text = "82÷7=11, 5"
$t.Cell(13, 5).Range.Text = "67÷7=9, 4"
$t.Cell(17, 1).Range.Text = "17÷5=3, 2"
$t.Cell(17, 2).Range.Text = "47÷9=5, 2"
$t.Cell(17, 3).Range.Text = "15÷5=3, 0"
$t.Cell(17, 4).Range.Text = "10÷6=1, 4"
$t.Cell(17, 5).Range.Text = "22÷3=7, 1"
